$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("I2:I17").NumberFormat = "@"
$wsResults.Range("N2:N17").NumberFormat = "@"

$rowData = @(1, "A", "Facility 1", "A", 100, 69900, "A", 0.5252, "1%", 0.519948, 363.443652, 699, 69536.556348, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(2, $c).Value = $rowData[$c-1] }
$rowData = @(1, "B", "Facility 1", "A", 100, 100, "B", 10, "1%", 9.7, 9.7, 1, 90.3, "5%", 0.485)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(3, $c).Value = $rowData[$c-1] }
$rowData = @(2, "A", "Facility 1", "B", 156, 1404000, "B", 70, "3%", 67.89999999999999, 611099.9999999999, 9000, 792900.0000000001, "5%", 30555)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(4, $c).Value = $rowData[$c-1] }
$rowData = @(3, "A", "Facility 4", "C", 423, 195849, "B", 65, "3%", 63.05, 29192.15, 463, 166656.85, "5%", 1459.6075)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(5, $c).Value = $rowData[$c-1] }
$rowData = @(3, "B", "Facility 4", "C", 423, 57951, "A", 55, "1%", 54.45, 7459.650000000001, 137, 50491.35, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(6, $c).Value = $rowData[$c-1] }
$rowData = @(4, "A", "Facility 4", "C", 453, 2447559, "A", 23, "1%", 22.77, 123026.31, 5403, 2324532.69, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(7, $c).Value = $rowData[$c-1] }
$rowData = @(4, "B", "Facility 4", "C", 453, 120951, "C", 24, "4%", 23.04, 6151.679999999999, 267, 114799.32, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(8, $c).Value = $rowData[$c-1] }
$rowData = @(5, "A", "Facility 5", "C", 342, 15390, "B", 34, "3%", 32.98, 1484.1, 45, 13905.9, "5%", 74.205)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(9, $c).Value = $rowData[$c-1] }
$rowData = @(6, "A", "Facility 9", "C", 653, 158026, "B", 24, "3%", 23.28, 5633.76, 242, 152392.24, "5%", 281.688)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(10, $c).Value = $rowData[$c-1] }
$rowData = @(7, "A", "Facility 9", "C", 432, 286848, "A", 23, "1%", 22.77, 15119.28, 664, 271728.72, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(11, $c).Value = $rowData[$c-1] }
$rowData = @(8, "A", "Facility 9", "C", 456, 10488, "B", 13, "3%", 12.61, 290.03, 23, 10197.97, "5%", 14.5015)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(12, $c).Value = $rowData[$c-1] }
$rowData = @(8, "B", "Facility 9", "C", 456, 456, "C", 24, "4%", 23.04, 23.04, 1, 432.96, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(13, $c).Value = $rowData[$c-1] }
$rowData = @(9, "A", "Facility 10", "C", 234, 54054, "C", 32, "4%", 30.72, 7096.32, 231, 46957.68, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(14, $c).Value = $rowData[$c-1] }
$rowData = @(9, "B", "Facility 10", "C", 234, 234, "B", 56, "3%", 54.32, 54.32, 1, 179.68, "5%", 2.716)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(15, $c).Value = $rowData[$c-1] }
$rowData = @(10, "A", "Facility 10", "C", 231, 2772, "B", 13, "3%", 12.61, 151.32, 12, 2620.68, "5%", 7.566)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(16, $c).Value = $rowData[$c-1] }
$rowData = @(10, "B", "Facility 10", "C", 231, 231, "C", 15, "4%", 14.4, 14.4, 1, 216.6, "0%", 0)
for ($c = 1; $c -le 15; $c++) { $wsResults.Cells.Item(17, $c).Value = $rowData[$c-1] }

$lpText = @"
\* Sourcing_with_MultiTier_Rebates_Discounts *\
Minimize
OBJ: S_A + S_B + S_C - rebate_A - rebate_B - rebate_C
Subject To
BaseSpend_A: S0_A - 0.5252 x_A_1 - 64 x_A_10 - 70 x_A_2 - 55 x_A_3 - 23 x_A_4
 - 54 x_A_5 - 42 x_A_6 - 23 x_A_7 - 75 x_A_8 - 97 x_A_9 = 0
BaseSpend_B: S0_B - 10 x_B_1 - 13 x_B_10 - 70 x_B_2 - 65 x_B_3 - 75 x_B_4
 - 34 x_B_5 - 24 x_B_6 - 85 x_B_7 - 13 x_B_8 - 56 x_B_9 = 0
BaseSpend_C: S0_C - 24 x_C_1 - 15 x_C_10 - 75 x_C_2 - 60 x_C_3 - 24 x_C_4
 - 44 x_C_6 - 42 x_C_7 - 24 x_C_8 - 32 x_C_9 = 0
Capacity_B_Bid_ID_1: x_B_1 <= 100000000
Capacity_B_Bid_ID_10: x_B_10 <= 100000000
Capacity_B_Bid_ID_2: x_B_2 <= 100000000
Capacity_B_Bid_ID_3: x_B_3 <= 100000000
Capacity_B_Bid_ID_4: x_B_4 <= 100000000
Capacity_B_Bid_ID_5: x_B_5 <= 100000000
Capacity_B_Bid_ID_6: x_B_6 <= 100000000
Capacity_B_Bid_ID_7: x_B_7 <= 100000000
Capacity_B_Bid_ID_8: x_B_8 <= 100000000
Capacity_B_Bid_ID_9: x_B_9 <= 100000000
Capacity_C_Bid_ID_1: x_C_1 <= 100000000
Capacity_C_Bid_ID_10: x_C_10 <= 100000000
Capacity_C_Bid_ID_2: x_C_2 <= 100000000
Capacity_C_Bid_ID_3: x_C_3 <= 100000000
Capacity_C_Bid_ID_4: x_C_4 <= 100000000
Capacity_C_Bid_ID_5: x_C_5 <= 100000000
Capacity_C_Bid_ID_6: x_C_6 <= 100000000
Capacity_C_Bid_ID_7: x_C_7 <= 100000000
Capacity_C_Bid_ID_8: x_C_8 <= 100000000
Capacity_C_Bid_ID_9: x_C_9 <= 100000000
Demand_1: x_A_1 + x_B_1 + x_C_1 = 700
Demand_10: x_A_10 + x_B_10 + x_C_10 = 13
Demand_2: x_A_2 + x_B_2 + x_C_2 = 9000
Demand_3: x_A_3 + x_B_3 + x_C_3 = 600
Demand_4: x_A_4 + x_B_4 + x_C_4 = 5670
Demand_5: x_A_5 + x_B_5 + x_C_5 = 45
Demand_6: x_A_6 + x_B_6 + x_C_6 = 242
Demand_7: x_A_7 + x_B_7 + x_C_7 = 664
Demand_8: x_A_8 + x_B_8 + x_C_8 = 24
Demand_9: x_A_9 + x_B_9 + x_C_9 = 232
DiscountTierLower_A_0: d_A - 19400000000 z_discount_A_0 >= -19400000000
DiscountTierLower_A_1: - 0.01 S0_A + d_A - 19400000000 z_discount_A_1
 >= -19400000000
DiscountTierLower_B_0: d_B - 97000000000 z_discount_B_0 >= -97000000000
DiscountTierLower_B_1: - 0.03 S0_B + d_B - 97000000000 z_discount_B_1
 >= -97000000000
DiscountTierLower_C_0: d_C - 97000000000 z_discount_C_0 >= -97000000000
DiscountTierLower_C_1: - 0.04 S0_C + d_C - 97000000000 z_discount_C_1
 >= -97000000000
DiscountTierMax_A_0: 19400000000 z_discount_A_0 <= 19400001000
DiscountTierMax_B_0: 97000000000 z_discount_B_0 <= 97000000500
DiscountTierMax_C_0: 97000000000 z_discount_C_0 <= 97000000500
_dummy: __dummy = 0
DiscountTierMin_A_0: __dummy >= 0
DiscountTierMin_A_1: x_A_1 + x_A_10 + x_A_3 + x_A_4 + x_A_8 + x_A_9
 - 1000 z_discount_A_1 >= 0
DiscountTierMin_B_0: __dummy >= 0
DiscountTierMin_B_1: x_B_2 + x_B_5 + x_B_6 + x_B_7 - 500 z_discount_B_1 >= 0
DiscountTierMin_C_0: __dummy >= 0
DiscountTierMin_C_1: x_C_1 + x_C_10 + x_C_3 + x_C_4 + x_C_8 + x_C_9
 - 500 z_discount_C_1 >= 0
DiscountTierSelect_A: z_discount_A_0 + z_discount_A_1 = 1
DiscountTierSelect_B: z_discount_B_0 + z_discount_B_1 = 1
DiscountTierSelect_C: z_discount_C_0 + z_discount_C_1 = 1
DiscountTierUpper_A_0: d_A + 19400000000 z_discount_A_0 <= 19400000000
DiscountTierUpper_A_1: - 0.01 S0_A + d_A + 19400000000 z_discount_A_1
 <= 19400000000
DiscountTierUpper_B_0: d_B + 97000000000 z_discount_B_0 <= 97000000000
DiscountTierUpper_B_1: - 0.03 S0_B + d_B + 97000000000 z_discount_B_1
 <= 97000000000
DiscountTierUpper_C_0: d_C + 97000000000 z_discount_C_0 <= 97000000000
DiscountTierUpper_C_1: - 0.04 S0_C + d_C + 97000000000 z_discount_C_1
 <= 97000000000
EffectiveSpend_A: - S0_A + S_A + d_A = 0
EffectiveSpend_B: - S0_B + S_B + d_B = 0
EffectiveSpend_C: - S0_C + S_C + d_C = 0
NonBid_C_5: x_C_5 = 0
RebateTierLower_A_0: rebate_A - 19400000000 y_rebate_A_0 >= -19400000000
RebateTierLower_A_1: - 0.1 S_A + rebate_A - 19400000000 y_rebate_A_1
 >= -19400000000
RebateTierLower_B_0: rebate_B - 97000000000 y_rebate_B_0 >= -97000000000
RebateTierLower_B_1: - 0.05 S_B + rebate_B - 97000000000 y_rebate_B_1
 >= -97000000000
RebateTierLower_C_0: rebate_C - 97000000000 y_rebate_C_0 >= -97000000000
RebateTierLower_C_1: - 0.07 S_C + rebate_C - 97000000000 y_rebate_C_1
 >= -97000000000
RebateTierMax_A_0: 19400000000 y_rebate_A_0 <= 19400000500
RebateTierMax_B_0: 97000000000 y_rebate_B_0 <= 97000000500
RebateTierMax_C_0: 97000000000 y_rebate_C_0 <= 97000000700
RebateTierMin_A_0: __dummy >= 0
RebateTierMin_A_1: - 500 y_rebate_A_1 >= 0
RebateTierMin_B_0: __dummy >= 0
RebateTierMin_B_1: x_B_1 + x_B_10 + x_B_3 + x_B_4 + x_B_8 + x_B_9
 - 500 y_rebate_B_1 >= 0
RebateTierMin_C_0: __dummy >= 0
RebateTierMin_C_1: - 700 y_rebate_C_1 >= 0
RebateTierSelect_A: y_rebate_A_0 + y_rebate_A_1 = 1
RebateTierSelect_B: y_rebate_B_0 + y_rebate_B_1 = 1
RebateTierSelect_C: y_rebate_C_0 + y_rebate_C_1 = 1
RebateTierUpper_A_0: rebate_A + 19400000000 y_rebate_A_0 <= 19400000000
RebateTierUpper_A_1: - 0.1 S_A + rebate_A + 19400000000 y_rebate_A_1
 <= 19400000000
RebateTierUpper_B_0: rebate_B + 97000000000 y_rebate_B_0 <= 97000000000
RebateTierUpper_B_1: - 0.05 S_B + rebate_B + 97000000000 y_rebate_B_1
 <= 97000000000
RebateTierUpper_C_0: rebate_C + 97000000000 y_rebate_C_0 <= 97000000000
RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1
 <= 97000000000
SupplIndicatorLB_0_A_1: - w_0_A_1 + x_A_1 >= 0
SupplIndicatorLB_0_A_10: - w_0_A_10 + x_A_10 >= 0
SupplIndicatorLB_0_A_3: - w_0_A_3 + x_A_3 >= 0
SupplIndicatorLB_0_A_4: - w_0_A_4 + x_A_4 >= 0
SupplIndicatorLB_0_A_8: - w_0_A_8 + x_A_8 >= 0
SupplIndicatorLB_0_A_9: - w_0_A_9 + x_A_9 >= 0
SupplIndicatorLB_0_B_1: - w_0_B_1 + x_B_1 >= 0
SupplIndicatorLB_0_B_10: - w_0_B_10 + x_B_10 >= 0
SupplIndicatorLB_0_B_3: - w_0_B_3 + x_B_3 >= 0
SupplIndicatorLB_0_B_4: - w_0_B_4 + x_B_4 >= 0
SupplIndicatorLB_0_B_8: - w_0_B_8 + x_B_8 >= 0
SupplIndicatorLB_0_B_9: - w_0_B_9 + x_B_9 >= 0
SupplIndicatorLB_0_C_1: - w_0_C_1 + x_C_1 >= 0
SupplIndicatorLB_0_C_10: - w_0_C_10 + x_C_10 >= 0
SupplIndicatorLB_0_C_3: - w_0_C_3 + x_C_3 >= 0
SupplIndicatorLB_0_C_4: - w_0_C_4 + x_C_4 >= 0
SupplIndicatorLB_0_C_8: - w_0_C_8 + x_C_8 >= 0
SupplIndicatorLB_0_C_9: - w_0_C_9 + x_C_9 >= 0
SupplIndicator_0_A_1: - 1000000000 w_0_A_1 + x_A_1 <= 0
SupplIndicator_0_A_10: - 1000000000 w_0_A_10 + x_A_10 <= 0
SupplIndicator_0_A_3: - 1000000000 w_0_A_3 + x_A_3 <= 0
SupplIndicator_0_A_4: - 1000000000 w_0_A_4 + x_A_4 <= 0
SupplIndicator_0_A_8: - 1000000000 w_0_A_8 + x_A_8 <= 0
SupplIndicator_0_A_9: - 1000000000 w_0_A_9 + x_A_9 <= 0
SupplIndicator_0_B_1: - 1000000000 w_0_B_1 + x_B_1 <= 0
SupplIndicator_0_B_10: - 1000000000 w_0_B_10 + x_B_10 <= 0
SupplIndicator_0_B_3: - 1000000000 w_0_B_3 + x_B_3 <= 0
SupplIndicator_0_B_4: - 1000000000 w_0_B_4 + x_B_4 <= 0
SupplIndicator_0_B_8: - 1000000000 w_0_B_8 + x_B_8 <= 0
SupplIndicator_0_B_9: - 1000000000 w_0_B_9 + x_B_9 <= 0
SupplIndicator_0_C_1: - 1000000000 w_0_C_1 + x_C_1 <= 0
SupplIndicator_0_C_10: - 1000000000 w_0_C_10 + x_C_10 <= 0
SupplIndicator_0_C_3: - 1000000000 w_0_C_3 + x_C_3 <= 0
SupplIndicator_0_C_4: - 1000000000 w_0_C_4 + x_C_4 <= 0
SupplIndicator_0_C_8: - 1000000000 w_0_C_8 + x_C_8 <= 0
SupplIndicator_0_C_9: - 1000000000 w_0_C_9 + x_C_9 <= 0
SupplierCount_0_1: w_0_A_1 + w_0_B_1 + w_0_C_1 >= 2
SupplierCount_0_10: w_0_A_10 + w_0_B_10 + w_0_C_10 >= 2
SupplierCount_0_3: w_0_A_3 + w_0_B_3 + w_0_C_3 >= 2
SupplierCount_0_4: w_0_A_4 + w_0_B_4 + w_0_C_4 >= 2
SupplierCount_0_8: w_0_A_8 + w_0_B_8 + w_0_C_8 >= 2
SupplierCount_0_9: w_0_A_9 + w_0_B_9 + w_0_C_9 >= 2
Transition_10_A: - 13 T_10_A + x_A_10 <= 0
Transition_10_B: - 13 T_10_B + x_B_10 <= 0
Transition_1_B: - 700 T_1_B + x_B_1 <= 0
Transition_1_C: - 700 T_1_C + x_C_1 <= 0
Transition_2_A: - 9000 T_2_A + x_A_2 <= 0
Transition_2_C: - 9000 T_2_C + x_C_2 <= 0
Transition_3_A: - 600 T_3_A + x_A_3 <= 0
Transition_3_B: - 600 T_3_B + x_B_3 <= 0
Transition_4_A: - 5670 T_4_A + x_A_4 <= 0
Transition_4_B: - 5670 T_4_B + x_B_4 <= 0
Transition_5_A: - 45 T_5_A + x_A_5 <= 0
Transition_5_B: - 45 T_5_B + x_B_5 <= 0
Transition_6_A: - 242 T_6_A + x_A_6 <= 0
Transition_6_B: - 242 T_6_B + x_B_6 <= 0
Transition_7_A: - 664 T_7_A + x_A_7 <= 0
Transition_7_B: - 664 T_7_B + x_B_7 <= 0
Transition_8_A: - 24 T_8_A + x_A_8 <= 0
Transition_8_B: - 24 T_8_B + x_B_8 <= 0
Transition_9_A: - 232 T_9_A + x_A_9 <= 0
Transition_9_B: - 232 T_9_B + x_B_9 <= 0
Volume_A: V_A - x_A_1 - x_A_10 - x_A_2 - x_A_3 - x_A_4 - x_A_5 - x_A_6 - x_A_7
 - x_A_8 - x_A_9 = 0
Volume_B: V_B - x_B_1 - x_B_10 - x_B_2 - x_B_3 - x_B_4 - x_B_5 - x_B_6 - x_B_7
 - x_B_8 - x_B_9 = 0
Volume_C: V_C - x_C_1 - x_C_10 - x_C_2 - x_C_3 - x_C_4 - x_C_5 - x_C_6 - x_C_7
 - x_C_8 - x_C_9 = 0
Binaries
T_10_A
T_10_B
T_1_B
T_1_C
T_2_A
T_2_C
T_3_A
T_3_B
T_4_A
T_4_B
T_5_A
T_5_B
T_6_A
T_6_B
T_7_A
T_7_B
T_8_A
T_8_B
T_9_A
T_9_B
w_0_A_1
w_0_A_10
w_0_A_3
w_0_A_4
w_0_A_8
w_0_A_9
w_0_B_1
w_0_B_10
w_0_B_3
w_0_B_4
w_0_B_8
w_0_B_9
w_0_C_1
w_0_C_10
w_0_C_3
w_0_C_4
w_0_C_8
w_0_C_9
y_rebate_A_0
y_rebate_A_1
y_rebate_B_0
y_rebate_B_1
y_rebate_C_0
y_rebate_C_1
z_discount_A_0
z_discount_A_1
z_discount_B_0
z_discount_B_1
z_discount_C_0
z_discount_C_1
End

"@
$wsLP = $wb.Worksheets.Item("LP Model")
$wsLP.Range("A2").Value = $lpText
